$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray C2 value entirely (naive forecaster bug produced an
# extra y_1 value for 2008 that should not have been written).
$ws.Range("C2").ClearContents()

# Refresh the recomputed forecast values (tiny precision corrections from
# fixing the naive component forecaster bug).
$ws.Range("E2").Value = 0.1825419310453436
$ws.Range("C5").Value = 1.692932643509826
$ws.Range("C6").Value = 1.020829760720687
$ws.Range("E7").Value = 1.258913537332873
$ws.Range("C8").Value = 1.019715257608933
$ws.Range("C9").Value = 2.173959184500363
$ws.Range("E9").Value = 1.566646323486065
$ws.Range("C10").Value = 1.707434489470039
$ws.Range("C11").Value = 1.456988786619817
$ws.Range("E11").Value = 1.84279714442821
$ws.Range("E12").Value = 1.580042106786372
$ws.Range("E13").Value = 1.604795846351492
$ws.Range("E14").Value = -1.215549235925817
$ws.Range("C15").Value = -5.665308402785508
$ws.Range("E15").Value = -4.458023117238186
$ws.Range("E16").Value = -0.5534294478199198
$ws.Range("C17").Value = -1.17492083522599
$ws.Range("E17").Value = 0.1140263184959744
$ws.Range("C18").Value = 0.02017133142706573
$ws.Range("E18").Value = -0.3342090768663986
$ws.Range("E19").Value = -0.1189552196680155
